# Generate Report for Handoff
# Update the status row for "ed93cd1b-bb04-41d1-b3e9-58feee8960e9.md" to reflect
# that it is now ready for handoff (was "In Translation").

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 corresponds to ed93cd1b-... .md ---
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-16-19 20:16:12"

# --- zh-cn sheet: row 3 corresponds to ed93cd1b-... .md ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "ed93cd1b-bb04-41d1-b3e9-58feee8960e9.91b95d3d80266ee97dfb02bb11fe732d85eda8a1.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-19 20:16:09"

# --- de-de sheet: row 3 corresponds to ed93cd1b-... .md ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "ed93cd1b-bb04-41d1-b3e9-58feee8960e9.91b95d3d80266ee97dfb02bb11fe732d85eda8a1.de-de.xlf"
$dede.Range("E3").Value = "2016-03-19 20:16:12"
